$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update CasesTab query cell (B2): weight formula now wrapped in CASE expression
$ws.Range("B2").Value() = 'MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
Match (c)<--(diag:diagnosis)
WHERE diag.disease_term in[''Pulmonary Neoplasms'']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co, demo.patient_age_at_enrollment AS age
RETURN  
       coalesce(c.case_id, '''') AS `Case ID`,
       coalesce(s.clinical_study_designation, '''') AS `Study Code`,
       coalesce(s.clinical_study_type, '''') AS  `Study Type`,
       coalesce(demo.breed, '''') AS Breed ,
       coalesce(diag.disease_term, '''') AS Diagnosis ,
       coalesce(diag.stage_of_disease, '''') AS `Stage of Disease`,
       CASE age % 1 WHEN 0 THEN apoc.convert.toInteger(age) ELSE age END AS Age,
       coalesce(demo.sex, '''') AS Sex,
       coalesce(demo.neutered_indicator, '''') AS `Neutered Status`,
       coalesce(CASE weight % 1 WHEN 0 THEN apoc.convert.toInteger(weight) ELSE weight END, '''') AS `Weight (kg)`,
       coalesce(diag.best_response, '''') AS `Response to Treatment`,
       coalesce(co.cohort_description, '''') AS `Cohort`'

# Update FilesTab query cell (B4): samp -> samp:sample type annotation added
$ws.Range("B4").Value() = 'MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f)-[*]->(samp:sample)
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE diag.disease_term IN [''Pulmonary Neoplasms'']

OPTIONAL MATCH (s:study)<--(c)<--(diag:diagnosis)<-[*]-(samp:sample)
WITH
        f, parent, c, demo, diag, s, samp,
        [''Bytes'', ''KB'', ''MB'', ''GB'', ''TB''] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, parent, c, demo, diag, s, samp,
        f.file_size /(1024^i) AS value, 
        10^precision AS factor,
        units[i] as unit
WITH    
        f, parent, c, demo, diag, s, samp, unit,
        round(factor * value)/factor AS size
RETURN 
        coalesce(f.file_name, '''') AS `File Name`,
        coalesce(f.file_type, '''') AS `File Type`,
        coalesce(labels(parent)[0], '''') AS `Association`,
        coalesce(f.file_description, '''') AS `Description`,
        coalesce(f.file_format, '''') AS `Format`,
        CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+'' '' +unit ELSE size+'' '' +unit END AS Size,
        coalesce(samp.sample_id, '''') AS `Sample ID`,
        coalesce(c.case_id, '''') AS `Case ID`,
        coalesce(demo.breed,'''') AS Breed ,
        coalesce(diag.disease_term,'''') AS Diagnosis'

# Row heights adjusted (minor layout recalculation)
$ws.Rows.Item(2).RowHeight() = 216.95
$ws.Rows.Item(3).RowHeight() = 225
$ws.Rows.Item(4).RowHeight() = 264
$ws.Rows.Item(5).RowHeight() = 197.1

# Column widths adjusted slightly (minor layout recalculation); column A left untouched to preserve bestFit
$ws.Columns.Item(2).ColumnWidth() = 86.877604166667
$ws.Columns.Item(3).ColumnWidth() = 74.877604166667
$ws.Columns.Item(4).ColumnWidth() = 69.451822916667
$ws.Columns.Item(5).ColumnWidth() = 27.736979166667
